$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the old row 4 (admin / password) - clearing makes the row disappear
# from the sheetData since it becomes fully empty, while row 5 keeps its own
# row index.
$ws.Range("A4:B4").ClearContents()

# Append the new rows at the bottom of the table.
$ws.Range("A6").Value = "Admin"
$ws.Range("B6").Value = "Qedge123!@#"

$ws.Range("A7").Value = "admin"
$ws.Range("B7").Value = "password"

# Update the active selection like Excel would persist on save.
$ws.Range("L11").Select()
